$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2: weight 65 -> 79
$ws.Range("E2").Value = 79

# K2: clear dob (was "2025-03-05") -> empty text
# Leading "'" forces text interpretation for the (empty) value; resetting the
# style back to Normal afterwards avoids leaving a stray quote-prefix format.
$ws.Range("K2").Value = "'"
$ws.Range("K2").Style = "Normal"

# R2: clear medical_conditions (was "pregnant") -> empty text
$ws.Range("R2").Value = "'"
$ws.Range("R2").Style = "Normal"

# S2: clear next_of_kin (was "u") -> empty text
$ws.Range("S2").Value = "'"
$ws.Range("S2").Style = "Normal"

# A3: id 3 (number) -> "2" (text)
$ws.Range("A3").Value = "'2"
$ws.Range("A3").Style = "Normal"

# C3: address "kalam" -> "jhoopri"
$ws.Range("C3").Value = "jhoopri"

# H3: payment_status "Paid" -> "Pending"
$ws.Range("H3").Value = "Pending"

# K3: dob "2025-03-01" -> "2025-03-03" (force text, not a date)
$ws.Range("K3").Value = "'2025-03-03"
$ws.Range("K3").Style = "Normal"

# R3: medical_conditions "strong AF" -> "strong"
$ws.Range("R3").Value = "strong"
